# Corrected return line type in sequence diagram
#
# - Re-draws (re-orders to front of z-order) the 5 lifeline "guide" connectors
#   so they sit right after the group shape properties, ahead of the rest of
#   the diagram's shapes (matches a delete+redraw of the lifelines).
# - Two of those lifeline connectors pick up new geometry/name (as if they
#   had been deleted and redrawn against their anchor rectangles).
# - A return-arrow connector gets switched from a solid line to the
#   System-Dash preset (the actual "corrected return line type").
# - A handful of labels/arrows shift slightly to line up with the corrected
#   lifelines.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# EMU -> Points helper. Shape.Left/Top/Width/Height are single-precision
# (points) under the hood, so nudge very slightly before the EMU->pt divide
# to avoid landing 1 EMU short after the round-trip.
function EMU($v) {
    return ($v + 0.6) / 12700
}

function GetShape($name) {
    return $s.Shapes.Item($name)
}

# ---------------------------------------------------------------------
# 1) Re-order the 5 lifeline connectors to the front of the z-order
#    (i.e. right after the spTree's own grpSpPr, ahead of "Group 96").
#    Send-to-back each one, starting with the one that should end up
#    furthest from the front, so the final relative order is:
#      Straight Connector 3, Straight Connector 11, Straight Connector 109,
#      Straight Connector 73, Straight Connector 22
# ---------------------------------------------------------------------
$cTmp = GetShape("Straight Connector 22")
$cTmp.ZOrder(1)
$cTmp = GetShape("Straight Connector 73")
$cTmp.ZOrder(1)
$cTmp = GetShape("Straight Connector 109")
$cTmp.ZOrder(1)
$cTmp = GetShape("Straight Connector 11")
$cTmp.ZOrder(1)
$cTmp = GetShape("Straight Connector 3")
$cTmp.ZOrder(1)

# ---------------------------------------------------------------------
# 2) Lifeline connectors that got redrawn with new geometry (and picked
#    up new shape names in the process).
# ---------------------------------------------------------------------
$connStorage = GetShape("Straight Connector 73")
$connStorage.Left   = EMU(10931636)
$connStorage.Top    = EMU(954530)
$connStorage.Width  = EMU(3064)
$connStorage.Height = EMU(4559506)
$connStorage.Name   = "Straight Connector 127"

$connCommand = GetShape("Straight Connector 22")
$connCommand.Left   = EMU(9625070)
$connCommand.Top    = EMU(2369481)
$connCommand.Width  = EMU(0)
$connCommand.Height = EMU(3172410)
$connCommand.Name   = "Straight Connector 133"

# ---------------------------------------------------------------------
# 3) The actual "corrected return line type": switch the return arrow to
#    the System Dash preset line style.
# ---------------------------------------------------------------------
$returnArrow = GetShape("Straight Arrow Connector 14")
$returnArrow.Line.DashStyle = 9   # msoLineLongDashDotDot -> prstDash "sysDash"

# ---------------------------------------------------------------------
# 4) Position/size touch-ups on labels & arrows to match the corrected
#    lifelines.
# ---------------------------------------------------------------------

# "executeCommand(...)" textbox
$shp = GetShape("TextBox 8")
$shp.Top = EMU(996024)

# ":Command" rectangle narrower
$shp = GetShape("Rectangle 21")
$shp.Width = EMU(1358853)

# "handleCollate(Command)" caption rectangle moves up slightly
$shp = GetShape("Rectangle 58")
$shp.Top = EMU(1863280)

# small activation-bar rectangle shrinks/moves
$shp = GetShape("Rectangle 63")
$shp.Left   = EMU(9501867)
$shp.Top    = EMU(2367513)
$shp.Height = EMU(125808)

# return arrow shifts up
$shp = GetShape("Straight Arrow Connector 64")
$shp.Top = EMU(2481773)

# "handleCollate" caption textbox
$shp = GetShape("TextBox 71")
$shp.Top = EMU(3094293)

# ":Storage" rectangle narrower & shifted
$shp = GetShape("Rectangle 72")
$shp.Left  = EMU(10347546)
$shp.Width = EMU(1168180)

# small activation-bar rectangle near Storage moves
$shp = GetShape("Rectangle 74")
$shp.Left = EMU(10802327)
$shp.Top  = EMU(4021257)

# long horizontal arrow shortens
$shp = GetShape("Straight Arrow Connector 75")
$shp.Width = EMU(6957544)

# activation-bar rectangle shortens
$shp = GetShape("Rectangle 77")
$shp.Width = EMU(8479703)

# bottom-most long dashed return arrow shortens
$shp = GetShape("Straight Arrow Connector 94")
$shp.Left  = EMU(3833277)
$shp.Width = EMU(6969050)
